$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DATA")
$ws.Activate()

# Move column O (remarks) data to column T, leaving P:S empty in between.
# Use Formula (not Value/Cut) for the copy so literal #N/A error cells keep
# their error type instead of collapsing to a raw VT_ERROR numeric code.
$src = $ws.Range("O1:O319")
$dst = $ws.Range("T1:T319")
$dst.Formula = $src.Formula()
$src.ClearContents()

# New blank columns O:S inherit a width matching column N (15.5 chars),
# which is what the sheet used for the moved "remarks" column.
$ws.Range("O1:S1").ColumnWidth = 15.5

# Update the hidden filter-database defined name so it spans the new
# right-hand edge of the data range (was $A$1:$O$319).
$name = $wb.Names.Item("DATA!_FilterDatabase")
$name.RefersTo = "=DATA!`$A`$1:`$T`$319"

# Restore the active selection on the frozen pane to B21.
$ws.Range("B21").Select()
